$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (incl. date number format/style) from row 5 down into the
# two new rows before filling in their values, so the new date cells (col G)
# keep the same style index as the existing date column.
$ws.Range("A5:H5").Copy($ws.Range("A6:H6"))
$ws.Range("A5:H5").Copy($ws.Range("A7:H7"))

# Row 6
$ws.Range("A6").Value = 9399.64
$ws.Range("B6").Value = 10438.24
$ws.Range("C6").Value = 22.31
$ws.Range("D6").Value = 20.09
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -9.9499999999999993
$ws.Range("G6").Value = 42607.884270833332
$ws.Range("H6").Value = $false

# Row 7
$ws.Range("A7").Value = 8880.7800000000007
$ws.Range("B7").Value = 9399.64
$ws.Range("C7").Value = 20.3
$ws.Range("D7").Value = 19.18
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = -5.52
$ws.Range("G7").Value = 42608.616354166668
$ws.Range("H7").Value = $false
